# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 13:22"

# --- Simple in-place numeric refreshes (no re-sort needed) ---

# Row 17: Brasil
$ws.Cells.Item(17, 2).Value = 16238
$ws.Cells.Item(17, 3).Value = 50
$ws.Cells.Item(17, 4).Value = 173
$ws.Cells.Item(17, 5).Value = 15242
$ws.Cells.Item(17, 7).Value = 3
$ws.Cells.Item(17, 8).Value = 823

# Row 19: Austria
$ws.Cells.Item(19, 2).Value = 13074
$ws.Cells.Item(19, 3).Value = 132
$ws.Cells.Item(19, 5).Value = 7539

# Row 75: Kazajistan
$ws.Cells.Item(75, 4).Value = 58
$ws.Cells.Item(75, 5).Value = 699

# --- India moves above Australia/Irlanda/Noruega (rows 24-27 re-sorted) ---

$ws.Cells.Item(24, 1).Value = "India"
$ws.Cells.Item(24, 2).Value = 6237
$ws.Cells.Item(24, 3).Value = 321
$ws.Cells.Item(24, 4).Value = 569
$ws.Cells.Item(24, 5).Value = 5482
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 8
$ws.Cells.Item(24, 8).Value = 186

$ws.Cells.Item(25, 1).Value = "Australia"
$ws.Cells.Item(25, 2).Value = 6104
$ws.Cells.Item(25, 3).Value = 52
$ws.Cells.Item(25, 4).Value = 2987
$ws.Cells.Item(25, 5).Value = 3066
$ws.Cells.Item(25, 6).Value = 81
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(25, 8).Value = 51

$ws.Cells.Item(26, 1).Value = "Irlanda"
$ws.Cells.Item(26, 2).Value = 6074
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 25
$ws.Cells.Item(26, 5).Value = 5814
$ws.Cells.Item(26, 6).Value = 165
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(26, 8).Value = 235

$ws.Cells.Item(27, 1).Value = "Noruega"
$ws.Cells.Item(27, 2).Value = 6042
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = 32
$ws.Cells.Item(27, 5).Value = 5909
$ws.Cells.Item(27, 6).Value = 78
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(27, 8).Value = 101

# --- Banglades moves above Ghana/San Marino/... (rows 101-111 re-sorted) ---

$ws.Cells.Item(101, 1).Value = "Banglades"
$ws.Cells.Item(101, 2).Value = 330
$ws.Cells.Item(101, 3).Value = 112
$ws.Cells.Item(101, 4).Value = 33
$ws.Cells.Item(101, 5).Value = 276
$ws.Cells.Item(101, 6).Value = 1
$ws.Cells.Item(101, 7).Value = 1
$ws.Cells.Item(101, 8).Value = 21

$ws.Cells.Item(102, 1).Value = "Ghana"
$ws.Cells.Item(102, 2).Value = 313
$ws.Cells.Item(102, 3).Value = 0
$ws.Cells.Item(102, 4).Value = 34
$ws.Cells.Item(102, 5).Value = 273
$ws.Cells.Item(102, 6).Value = 2
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 6

$ws.Cells.Item(103, 1).Value = "San Marino"
$ws.Cells.Item(103, 2).Value = 308
$ws.Cells.Item(103, 3).Value = 0
$ws.Cells.Item(103, 4).Value = 45
$ws.Cells.Item(103, 5).Value = 229
$ws.Cells.Item(103, 6).Value = 14
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 34

$ws.Cells.Item(104, 1).Value = "Kirguistan"
$ws.Cells.Item(104, 2).Value = 280
$ws.Cells.Item(104, 3).Value = 10
$ws.Cells.Item(104, 4).Value = 35
$ws.Cells.Item(104, 5).Value = 241
$ws.Cells.Item(104, 6).Value = 5
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 4

$ws.Cells.Item(105, 1).Value = "Nigeria"
$ws.Cells.Item(105, 2).Value = 276
$ws.Cells.Item(105, 3).Value = 0
$ws.Cells.Item(105, 4).Value = 44
$ws.Cells.Item(105, 5).Value = 226
$ws.Cells.Item(105, 6).Value = 2
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = 6

$ws.Cells.Item(106, 1).Value = "Mauricio"
$ws.Cells.Item(106, 2).Value = 273
$ws.Cells.Item(106, 3).Value = 0
$ws.Cells.Item(106, 4).Value = 19
$ws.Cells.Item(106, 5).Value = 247
$ws.Cells.Item(106, 6).Value = 3
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 7

$ws.Cells.Item(107, 1).Value = "Bolivia"
$ws.Cells.Item(107, 2).Value = 264
$ws.Cells.Item(107, 3).Value = 54
$ws.Cells.Item(107, 4).Value = 2
$ws.Cells.Item(107, 5).Value = 244
$ws.Cells.Item(107, 6).Value = 3
$ws.Cells.Item(107, 7).Value = 3
$ws.Cells.Item(107, 8).Value = 18

$ws.Cells.Item(108, 1).Value = "Estado de Palestina"
$ws.Cells.Item(108, 2).Value = 263
$ws.Cells.Item(108, 3).Value = 0
$ws.Cells.Item(108, 4).Value = 44
$ws.Cells.Item(108, 5).Value = 218
$ws.Cells.Item(108, 6).Value = 0
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 1

$ws.Cells.Item(109, 1).Value = "Montenegro"
$ws.Cells.Item(109, 2).Value = 252
$ws.Cells.Item(109, 3).Value = 4
$ws.Cells.Item(109, 4).Value = 4
$ws.Cells.Item(109, 5).Value = 246
$ws.Cells.Item(109, 6).Value = 7
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 2

$ws.Cells.Item(110, 1).Value = "Vietnam"
$ws.Cells.Item(110, 2).Value = 251
$ws.Cells.Item(110, 3).Value = 0
$ws.Cells.Item(110, 4).Value = 128
$ws.Cells.Item(110, 5).Value = 123
$ws.Cells.Item(110, 6).Value = 8
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 0

$ws.Cells.Item(111, 1).Value = "Senegal"
$ws.Cells.Item(111, 2).Value = 244
$ws.Cells.Item(111, 3).Value = 0
$ws.Cells.Item(111, 4).Value = 113
$ws.Cells.Item(111, 5).Value = 129
$ws.Cells.Item(111, 6).Value = 1
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 2
